$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 10: close out the previous block - A10 gets the "closing" border style (like A3),
#     and B10:E10 become empty cells carrying the matching format only (like B3:E3). ---
$ws.Range("A3").Copy() | Out-Null
$ws.Range("A10").PasteSpecial(-4122) | Out-Null
$ws.Range("B3:E3").Copy() | Out-Null
$ws.Range("B10:E10").PasteSpecial(-4122) | Out-Null

# --- Rows 11-12: new "mid-block" entries, formatted like row 4 (A=8,B=8,C/D/E=9) ---
$ws.Range("A4:E4").Copy() | Out-Null
$ws.Range("A11:E12").PasteSpecial(-4122) | Out-Null

# --- Row 13: new block start, formatted like row 2 (A=4,B=4,C/D/E=5) ---
$ws.Range("A2:E2").Copy() | Out-Null
$ws.Range("A13").PasteSpecial(-4122) | Out-Null

# --- Rows 14-16: continuation rows (no filename in A), formatted like row 2 columns B:E ---
$ws.Range("B2:E2").Copy() | Out-Null
$ws.Range("B14:E16").PasteSpecial(-4122) | Out-Null

# --- Fill in the numbers (do not touch shared-string table) ---
$ws.Range("B11").Value = 235
$ws.Range("B12").Value = 216
$ws.Range("B13").Value = 197
$ws.Range("B14").Value = 175
$ws.Range("B15").Value = 178
$ws.Range("B16").Value = 156

# --- Fill in the text, in the exact order the strings were authored (preserves sharedStrings.xml order) ---
$ws.Range("A11").Value = "SCRIPT/T01P01A/us0107.ssb"
$ws.Range("C11").Value = " Thank you for saving the world!"
$ws.Range("D11").Value = " Спасибо вам за спасение мира!"
$ws.Range("E11").Value = " Òðàòéáï âàí èà òðàòåîéå íéñà!"
$ws.Range("C12").Value = " That mountain almost reaches\npast the sky! I can\'t wait to climb it! ♪"
$ws.Range("A12").Value = "SCRIPT/P01P04A/us3102.ssb"
$ws.Range("E12").Value = " Âåñšéîà üóïê ãïñú äïòóéãàåó\nîåáåò! Ÿ óàë öïœô îà îåæ âèïáñàóûòÿ! ♪"
$ws.Range("D12").Value = " Вершина этой горы достигает\nнебес! Я так хочу на неё взобраться! ♪"
$ws.Range("C13").Value = " Wow! ♪ A flower garden! ♪"
$ws.Range("A13").Value = "SCRIPT/D73P11A/us3122.ssb"
$ws.Range("E13").Value = " Âàô! ♪ Øâåóïœîúê òàä! ♪"
$ws.Range("D13").Value = " Вау! ♪ Цветочный сад! ♪"
$ws.Range("C14").Value = " I wonder who made the\nSecret Room."
$ws.Range("C15").Value = " It\'s so mysterious! ♪"
$ws.Range("D14").Value = " Интересно, кто создал Секретную\nКомнату."
$ws.Range("D15").Value = " Как таинственно! ♪"
$ws.Range("E14").Value = " Éîóåñåòîï, ëóï òïèäàì Òåëñåóîôý\nËïíîàóô."
$ws.Range("E15").Value = " Ëàë óàéîòóâåîîï! ♪"
$ws.Range("C16").Value = " Those [CS:N]Shaymin[CR] have a wide\nknowledge of many things! ♪"
$ws.Range("D16").Value = " [CS:N]Шеймины[CR] о многом знают! ♪"
$ws.Range("E16").Value = " [CS:N]Šåêíéîú[CR] ï íîïãïí èîàýó! ♪"

# --- Row heights (explicit, matching the authored sizes) ---
$ws.Rows.Item(11).RowHeight = 43.2
$ws.Rows.Item(12).RowHeight = 43.2
$ws.Rows.Item(13).RowHeight = 43.2
$ws.Rows.Item(14).RowHeight = 27.6
$ws.Rows.Item(16).RowHeight = 21.6

# --- View: scrolled down to row 13, D16 selected ---
$ws.Application.ActiveWindow.ScrollRow = 13
$ws.Range("D16").Select() | Out-Null

Write-Output "edit complete"
